# [Shubham] Push code with new functionality
# Target sheet is the active sheet ("Sheet6" in tab order, physically
# stored as xl/worksheets/sheet4.xml) which holds the "Test on Flight 1"
# login/search data row used by the automation suite.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet6")

# ---------------------------------------------------------------------
# 1. Refresh the credentials & dates on row 2 (existing OneWay test row)
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "ankur.yadav@quadlabs.com"
$ws.Range("H2").Copy()
$ws.Range("G2").PasteSpecial(-4122)              # xlPasteFormats - pick up H2's hyperlink-cell style
$ws.Range("H2").Value = "Ankur@123"

$ws.Range("O2").Value = "25-Nov-2024"
$ws.Range("P2").Value = "29-Nov-2024"

foreach ($h in $ws.Hyperlinks) {
    $r = $h.Range.Address()
    if ($r -eq "`$G`$2") { $h.Address = "mailto:ankur.yadav@quadlabs.com" }
    if ($r -eq "`$H`$2") { $h.Address = "mailto:Ankur@123" }
}

# ---------------------------------------------------------------------
# 2. Duplicate row 2 into a new row 3 (values + formats) for a second,
#    RoundTrip variant of the same test case.
# ---------------------------------------------------------------------
$ws.Range("A2:Y2").Copy($ws.Range("A3:Y3"))
$ws.Rows.Item(3).RowHeight = 24.75

$ws.Range("I3").Value = "RoundTrip"

$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:ankur.yadav@quadlabs.com")
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:Ankur@123")

# Hyperlinks.Add() resets the font to the default hyperlink style, so
# re-apply G2's/H2's already-corrected formatting on top of it.
$ws.Range("H2").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("H3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Extend the data-validation dropdowns so they also cover row 3.
# ---------------------------------------------------------------------
$pairs = @(
    @("W2","W2:W3","1,2,3,4,5,6,7,8,9,10"),
    @("S2","S2:S3","1,2,3,4,5,6,7,8,9,10"),
    @("Y2","Y2:Y3","Basic,Semi Flex,Flex"),
    @("U2","U2:U3","Basic,Semi Flex,Flex"),
    @("T2","T2:T3","Saver,Basic,Semi Flex,Flex"),
    @("X2","X2:X3","Saver,Basic,Semi Flex,Flex"),
    @("V2","V2:V3","Guest,Business"),
    @("R2","R2:R3","Guest,Business"),
    @("Q2","Q2:Q3","Employee,Guest,MultiPassenger"),
    @("J2","J2:J3","Domestic,International"),
    @("I2","I2:I3","OneWay,RoundTrip,MultiCity"),
    @("F2","F2:F3","Administrator,Travel Arranger,Employee"),
    @("E2","E2:E3","Saurabh,Prince Chaurasia,Gunjan Swain,Shubham,Laxmi Khanal,Sudesh Kumar,Piyush,Ankur,D Divaker S,Ankur Yadav,Sachin Kumar")
)

foreach ($p in $pairs) {
    $ws.Range($p[0]).Validation.Delete()
}
foreach ($p in $pairs) {
    $ws.Range($p[1]).Validation.Add(3, 1, 1, $p[2])
}

# ---------------------------------------------------------------------
# 4. Best-fit the columns (values got wider with the new row/content).
# ---------------------------------------------------------------------
$widths = @{
    1  = 13.592447916666666
    2  = 62.592447916666664
    3  = 9.166666666666666
    4  = 21.022135416666668
    5  = 12.022135416666666
    6  = 12.592447916666666
    7  = 26.022135416666668
    8  = 11.451822916666666
    9  = 7.877604166666667
    10 = 10.736979166666666
    11 = 13.736979166666666
    12 = 28.877604166666668
    13 = 18.736979166666668
    14 = 27.736979166666668
    15 = 13.451822916666666
    16 = 10.307291666666666
    17 = 9.022135416666666
    19 = 6.166666666666667
    20 = 5.307291666666667
    21 = 8.736979166666666
    22 = 9.451822916666666
    23 = 7.307291666666667
    24 = 6.451822916666667
    25 = 8.736979166666666
}
foreach ($c in $widths.Keys) {
    $ws.Columns.Item($c).ColumnWidth = $widths[$c]
}

# ---------------------------------------------------------------------
# 5. Misc view/print state.
# ---------------------------------------------------------------------
$ws.Range("H9").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
